# Apply the workbook edit:
#  1. Insert a new "Player Info" sheet before the existing "ODI Batting" sheet,
#     with player ID/NAME/BATTING_HAND/BOWL_STYLE data.
#  2. In "ODI Batting", rename MATCH_CARD_LINK -> MATCH_CODE and replace the
#     full scorecard URLs with just the numeric match code.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Player Info" sheet (Worksheets.Add() inserts it before
#     the currently active sheet, matching its position ahead of "ODI Batting"). ---
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Header row
$playerInfo.Cells.Item(1, 1).Value = "ID"
$playerInfo.Cells.Item(1, 2).Value = "NAME"
$playerInfo.Cells.Item(1, 3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1, 4).Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row (force text so numeric-looking values stay strings, like the source)
$playerInfo.Cells.Item(2, 1).Value = "'4278"
$playerInfo.Cells.Item(2, 2).Value = "Mark Sinclair Chapman"
$playerInfo.Cells.Item(2, 3).Value = "Left Handed"
$playerInfo.Cells.Item(2, 4).Value = "Left Arm Orthodox"

$playerInfo.Range("A1").Select()

# --- 2. Update "ODI Batting": header + MatchCode-only values for column D ---
# Re-fetch the sheet reference now (after the insert) so it resolves correctly.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$matchCodes = @{ 2 = "3860"; 3 = "3862"; 4 = "4138"; 5 = "4139"; 6 = "4149"; 7 = "4406"; 8 = "4625" }
foreach ($row in $matchCodes.Keys) {
    $battingSheet.Cells.Item($row, 4).Value = "'" + $matchCodes[$row]
}
